# Update the "Rules" worksheet cell E8 from "Good Morning" to "GIT UPDATE"
# and move the active selection to E8 (as reflected in the sheetView).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
